# Auto update: 2025-12-05 02:00:49
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date (column A) moves from 2025-12-03 -> 2025-12-05.
# Force text formatting first so Excel doesn't auto-convert the
# ISO-looking string into a date serial number, then drop back to the
# Normal style so no lingering number-format is left on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-05"
$ws.Range("A2").Style = "Normal"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-12-05"
$ws.Range("A3").Style = "Normal"

# MACRO_SIGNAL (column O) text updates
$ws.Range("O2").Value = "⚪ 중립 구간"
$ws.Range("O3").Value = "⚪ 중립 구간"

# Row 2 (Joby / JOBY) numeric updates
$ws.Range("D2").Value = 15.55
$ws.Range("E2").Value = 56.1
$ws.Range("F2").Value = 10.13
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 57.1
$ws.Range("N2").Value = 52.43913937059539

# Row 3 (Archer / ACHR) numeric updates
$ws.Range("D3").Value = 8.68
$ws.Range("E3").Value = 58.8
$ws.Range("F3").Value = 15.83
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 53.7
$ws.Range("N3").Value = 52.43913937059539
